# "add backend+frontend to download preview data"
#
# The preview sheet used to stop after a handful of partially-filled demo
# rows (11 rows, several columns blank because the old backend only sent
# a sparse sample, plus one stray "ashish" test record). The new backend
# now returns the full, densely-populated preview rows for every "friend"
# record - 23 data rows (A2:G24), each with every column (name/bank/accNo
# /addr/frndName/frndBest) populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues
$fmtFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# G2 / G8 already hold the literal text "false" / "true" (as plain shared
# strings, not booleans). Stash copies of them off to the side first -
# assigning the bare word via Range.Value later on would otherwise get
# auto-coerced to a real Boolean cell, which is not what the sheet wants.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial($fmtValues) | Out-Null
$ws.Range("G8").Copy() | Out-Null
$ws.Range("Z2").PasteSpecial($fmtValues) | Out-Null

# Each entry: row index (1-based, data starts at row 2), then the values for
# columns A..G. Strings are prefixed "s:" so the loop below writes them as
# text even when they look numeric (e.g. "null" in the accNo column). The
# frnds_best column uses the placeholders "FALSE"/"TRUE" (handled specially
# below) since literal "false"/"true" text can't be assigned through .Value.
$rows = @(
    @(2,  0,  "s:abhi",  "s:bank1", 2345,     "s:null",        "s:frnd1",  "FALSE"),
    @(3,  1,  "s:abhi",  "s:bank2", 3456,     "s:address1111", "s:frnd2",  "s:null"),
    @(4,  2,  "s:abhi",  "s:bank1", 2345,     "s:null",        "s:frnd3",  "s:null"),
    @(5,  3,  "s:abhi",  "s:bank2", 3456,     "s:address1111", "s:frnd4",  "s:null"),
    @(6,  4,  "s:abhi",  "s:bank1", 2345,     "s:null",        "s:frnd1",  "FALSE"),
    @(7,  5,  "s:abhi",  "s:bank2", 3456,     "s:address1111", "s:frnd2",  "s:null"),
    @(8,  6,  "s:abhi",  "s:bank1", 2345,     "s:null",        "s:frnd3",  "s:null"),
    @(9,  7,  "s:abhi",  "s:bank2", 3456,     "s:address1111", "s:frnd4",  "s:null"),
    @(10, 8,  "s:aditi", "s:bank3", 98765,    "s:null",        "s:frnd5",  "FALSE"),
    @(11, 9,  "s:aditi", "s:bank4", 12345,    "s:null",        "s:frnd6",  "s:null"),
    @(12, 10, "s:aditi", "s:bank5", "s:null", "s:address2222", "s:frnd8",  "TRUE"),
    @(13, 11, "s:aditi", "s:bank3", 98765,    "s:null",        "s:frnd9",  "s:null"),
    @(14, 12, "s:aditi", "s:bank4", 12345,    "s:null",        "s:frnd10", "s:null"),
    @(15, 13, "s:aditi", "s:bank5", "s:null", "s:address2222", "s:frnd5",  "FALSE"),
    @(16, 14, "s:aditi", "s:bank3", 98765,    "s:null",        "s:frnd6",  "s:null"),
    @(17, 15, "s:aditi", "s:bank4", 12345,    "s:null",        "s:frnd8",  "TRUE"),
    @(18, 16, "s:aditi", "s:bank5", "s:null", "s:address2222", "s:frnd9",  "s:null"),
    @(19, 17, "s:aditi", "s:bank3", 98765,    "s:null",        "s:frnd10", "s:null"),
    @(20, 18, "s:aditi", "s:bank4", 12345,    "s:null",        "s:frnd5",  "FALSE"),
    @(21, 19, "s:aditi", "s:bank5", "s:null", "s:address2222", "s:frnd6",  "s:null"),
    @(22, 20, "s:aditi", "s:bank3", 98765,    "s:null",        "s:frnd8",  "TRUE"),
    @(23, 21, "s:aditi", "s:bank4", 12345,    "s:null",        "s:frnd9",  "s:null"),
    @(24, 22, "s:aditi", "s:bank5", "s:null", "s:address2222", "s:frnd10", "s:null")
)

$cols = @("A", "B", "C", "D", "E", "F", "G")
$boolFixups = @()

foreach ($entry in $rows) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $entry[$i + 1]
        if ($val -eq "FALSE" -or $val -eq "TRUE") {
            # defer - fill in via the stashed literal-string cells below
            $boolFixups += , @("$col$r", $val)
            continue
        }
        if ($val -is [string] -and $val.StartsWith("s:")) {
            $val = $val.Substring(2)
        }
        $ws.Range("$col$r").Value = $val
    }
}

foreach ($fix in $boolFixups) {
    $addr = $fix[0]
    $which = $fix[1]
    if ($which -eq "FALSE") {
        $ws.Range("Z1").Copy() | Out-Null
    } else {
        $ws.Range("Z2").Copy() | Out-Null
    }
    $ws.Range($addr).PasteSpecial($fmtValues) | Out-Null
}

# drop the scratch cells used to stash the literal "false"/"true" text
$ws.Range("Z1:Z2").ClearContents()

# Carry the row-number column's border/centered style (same as the header
# row used for A2 originally) down across all the newly added rows.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2:A24").PasteSpecial($fmtFormats) | Out-Null
$excel.CutCopyMode = $false
